$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new row for "d=6" is being inserted right before the existing
# "d=7" row, so the d=7 row (old row 7) and the d=10 row (old row 8)
# each shift down by one. Shift the bottom two rows down first so we
# don't clobber data we still need to move.

# Give the about-to-be-created row 9 the same formatting as the
# existing label rows (bold / centered / bordered), then move the old
# row 8 ("d=10") values into it.
$ws.Range("A6").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A9").Value = $ws.Range("A8").Value2
$ws.Range("B9").Value = $ws.Range("B8").Value2
$ws.Range("C9").Value = $ws.Range("C8").Value2
$ws.Range("D9").Value = $ws.Range("D8").Value2
$ws.Range("E9").Value = $ws.Range("E8").Value2

# Move the old row 7 ("d=7") values down into row 8.
$ws.Range("A8").Value = $ws.Range("A7").Value2
$ws.Range("B8").Value = $ws.Range("B7").Value2
$ws.Range("C8").Value = $ws.Range("C7").Value2
$ws.Range("D8").Value = $ws.Range("D7").Value2
$ws.Range("E8").Value = $ws.Range("E7").Value2

# Write the new "d=6" row (Diebold-Mariano correction) into row 7.
$ws.Range("A7").Value = "d=6"
$ws.Range("B7").Value = 98.21046445653728
$ws.Range("C7").Value = 97.61420596324969
$ws.Range("D7").Value = 98.57091521349003
$ws.Range("E7").Value = 97.49885176329013
